$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Backlog")

# Update the shared-string text values in D23 / E23 to reflect the new
# "copy path+name directly when in the path column" behaviour.
$ws.Range("D23").Value = "to copy the path and name when we are in the path column"
$ws.Range("E23").Value = "to have directly the path and name"

# Column D needs to grow to fit the longer text (closest reachable width).
$ws.Columns.Item(4).ColumnWidth = 50

# Move the active selection from F23 to F24.
$ws.Range("F24").Select()
